$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT even when it looks numeric
# (Excel's COM layer auto-converts plain numeric-looking strings to
# numbers on assignment, so force the cell to Text format first, then
# restore the cell's normal style so no stray formatting is left behind).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Update the "code" column (B) values - group codes reshuffled for test cases
$ws.Range("B1").Value = "9X"
Set-TextValue $ws.Range("B3") "16"
Set-TextValue $ws.Range("B6") "97"
$ws.Range("B5").Value = "E2"
$ws.Range("B7").Value = "OVX"
$ws.Range("B8").Value = "0J"
$ws.Range("B9").Value = "DK"
$ws.Range("B12").Value = "5P"

# Select whole row 2 (as if the user clicked the row header for row 2),
# matching the new selection state.
$ws.Rows(2).Select()
